# Remove the "John", "Jane", "James", "Jack" rows (old rows 2-5) and the
# "Helem" row (old row 9), shifting the remaining data rows (Rose, Beth,
# Kate, Niko) up so they become the new rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Helem" row first (row 9) so the row numbers used for the
# second delete (rows 2-5) remain valid.
$ws.Rows("9:9").Delete()
$ws.Rows("2:5").Delete()

# Update the active selection to match the post-edit state (A7).
$ws.Range("A7").Select()
